$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old "Start" row (was row 2) by deleting it, shifting everything up.
$ws.Rows.Item(2).Delete()

# Update the labels (shared strings will be re-used/rebuilt by the engine)
$ws.Range("A2").Value = "Start Scrum 1"
$ws.Range("A3").Value = "End Scrum 1"
$ws.Range("A4").Value = "Start Scrum 2"
$ws.Range("A5").Value = "End Scrum 2"
$ws.Range("A6").Value = "Start Scrum 3"
$ws.Range("A7").Value = "End Scrum 3"

# Update the burndown values
$ws.Range("B2").Value = 23
$ws.Range("B3").Value = 17
$ws.Range("B4").Value = 7
$ws.Range("B5").Value = 5
$ws.Range("B6").Value = 0
$ws.Range("B7").Value = 0

# Update selection to match the committed state
$ws.Range("B6").Select()
